$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.793.52'
$ws.Range("E2").Value = '  -3.11%  '

$ws.Range("D3").Value = '3.778.74'
$ws.Range("E3").Value = '  +1.20%  '

$ws.Range("E4").Value = '  +0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.33'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.23%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.72'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.78%  '

$ws.Range("D7").Value = '3.773.87'
$ws.Range("E7").Value = '  +1.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.533'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.33%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.158'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.49%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.29'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.72%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.467'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.73%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '38.15'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.38%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000242'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.41%  '

$ws.Range("D15").Value = '4.411.26'
$ws.Range("E15").Value = '  +1.55%  '

$ws.Range("D16").Value = '3.783.04'
$ws.Range("E16").Value = '  +1.60%  '

$ws.Range("D17").Value = '67.939.34'
$ws.Range("E17").Value = '  -2.80%  '

$ws.Range("E18").Value = '  -4.52%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.21'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.56%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.96'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '485.11'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.28%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.31'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.12%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.732'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.78%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.63'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.03%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.36'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -7.21%  '

$ws.Range("E26").Value = '  +4.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.19'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.89%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.11'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -9.91%  '

$ws.Range("E29").Value = '  -0.09%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.92'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.26%  '

$ws.Range("E31").Value = '  -1.36%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '32.31'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.50%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.55'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.77%  '

$ws.Range("E34").Value = '  -3.76%  '

$ws.Range("E35").Value = '  +0.16%  '

$ws.Range("E36").Value = '  -4.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.81'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.36%  '

$ws.Range("E38").Value = '  -2.38%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.324'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.17%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '442.67'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.92%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '48.89'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.24%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.32%  '

$ws.Range("E43").Value = '  -8.32%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.30'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.20%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.28'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -6.74%  '

$ws.Range("D46").Value = '2.842.47'
$ws.Range("E46").Value = '  -3.69%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0351'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.38%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '137.29'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.64%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '26.24'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.31%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.31'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.05%  '
